$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared strings / text values
$ws.Range("B9").Value = "COMPLEX, STA. ROSA, LAGUNA"
$ws.Range("A10").Value = "COMPLEX TECH"

# Update the invoice date (H9)
$ws.Range("H9").Value = 45316.66666666667

# Update the transaction date (B16) and amount (H16) - placeholder for transaction number
$ws.Range("B16").Value = 45293.8780324074
$ws.Range("H16").Value = 123
